$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update baseline (B) values, clear DALI (C) values, update INR (D) values
# Row 2 (nw 1)
$ws.Range("B2").Value = 253.3
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 221

# Row 3 (nw 2)
$ws.Range("B3").Value = 132.5
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 221

# Row 4 (nw 3)
$ws.Range("B4").Value = 126.3
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 221

# Row 5 (nw 4)
$ws.Range("B5").Value = 128.2
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 221

# Row 6 (nw 5)
$ws.Range("B6").Value = 127.4
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 221

# Row 7 (nw 6)
$ws.Range("B7").Value = 128.2
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 221

# Row 8 (nw 7)
$ws.Range("B8").Value = 128.5
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 221

# Row 9 (nw 8)
$ws.Range("B9").Value = 129.3
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 221

# Update selection to C2:C9 with active cell C2 (matches the authored view state)
$ws.Range("C2:C9").Select()
